$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in CodiceFiscale (B) and NumeroEntiAttivi (C) are stored as text in the
# source data (values can have leading zeros), so force text format on exactly
# the cells being written before assigning their new values.
$textCells = @(
    'C3', 'C4', 'C5', 'C6', 'C7', 'C8', 'C9', 'C11', 'C13', 'C14', 'C15', 'C16', 'C17', 'C18', 'C19', 'C20',
    'C21', 'C22', 'B23', 'C23', 'B24', 'C24', 'B25', 'C25', 'C26', 'B27', 'C27', 'B28', 'C28', 'C29', 'C31',
    'C34', 'C35', 'C38', 'B39', 'C39', 'B40', 'C41', 'B42', 'B43', 'C43', 'C44', 'C45', 'B47', 'C47', 'B49',
    'B53', 'C53', 'B54', 'C54', 'B59', 'C59', 'B60', 'B61', 'C61', 'B62', 'C62', 'B63', 'B64', 'B65', 'C65',
    'B66', 'C66', 'B67', 'C67', 'B68', 'B69', 'B70', 'B71', 'B74', 'B75', 'B76', 'B77', 'B78', 'B79', 'C79'
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$changes = @(
    ,@('C3', '1078')
    ,@('C4', '802')
    ,@('C5', '585')
    ,@('C6', '502')
    ,@('C7', '436')
    ,@('C8', '424')
    ,@('C9', '415')
    ,@('C11', '350')
    ,@('C13', '247')
    ,@('C14', '239')
    ,@('C15', '219')
    ,@('C16', '197')
    ,@('C17', '183')
    ,@('C18', '174')
    ,@('C19', '168')
    ,@('C20', '163')
    ,@('C21', '155')
    ,@('C22', '124')
    ,@('A23', 'Regione Toscana')
    ,@('B23', '01386030488')
    ,@('C23', '110')
    ,@('A24', 'Next Step Solution')
    ,@('B24', '02554480349')
    ,@('C24', '109')
    ,@('A25', 'Regione Basilicata')
    ,@('B25', '80002950766')
    ,@('C25', '106')
    ,@('C26', '95')
    ,@('A27', 'Siscom SPA')
    ,@('B27', '01778000040')
    ,@('C27', '92')
    ,@('A28', 'ANCITEL')
    ,@('B28', '07196850585')
    ,@('C28', '86')
    ,@('C29', '76')
    ,@('C31', '67')
    ,@('C34', '54')
    ,@('C35', '50')
    ,@('C38', '43')
    ,@('A39', 'UNIMATICA S.P.A')
    ,@('B39', '02098391200')
    ,@('C39', '41')
    ,@('A40', 'Consorzio I.T. Srl')
    ,@('B40', '01321400192')
    ,@('C41', '37')
    ,@('A42', 'Regione Liguria')
    ,@('B42', '00849050109')
    ,@('A43', 'Unicredit, Societa'' per Azioni')
    ,@('B43', '00348170101')
    ,@('C43', '33')
    ,@('C44', '28')
    ,@('C45', '27')
    ,@('A47', 'Nexi SpA')
    ,@('B47', '13212880150')
    ,@('C47', '19')
    ,@('A49', 'Regione Lazio')
    ,@('B49', '80143490581')
    ,@('A53', 'Aric Agenzia Regionale di Informatica e Committenza')
    ,@('B53', '91022630676')
    ,@('C53', '13')
    ,@('A54', 'Crédit Agricole Group Solutions Società Consortile per azioni')
    ,@('B54', '02771790348')
    ,@('C54', '12')
    ,@('A59', 'Argentea S.r.l.')
    ,@('B59', '02260390220')
    ,@('C59', '5')
    ,@('A60', 'Phoenix IT Solutions S.r.L')
    ,@('B60', '07623321218')
    ,@('A61', 'ARGO SOFTWARE SRL')
    ,@('B61', '00838520880')
    ,@('C61', '4')
    ,@('A62', 'CityPoste Payment Digital S.r.l.')
    ,@('B62', '02003750672')
    ,@('C62', '4')
    ,@('A63', 'e-SED Società Cooperativa')
    ,@('B63', '02695640421')
    ,@('A64', 'Linea Comune Spa')
    ,@('B64', '05591950489')
    ,@('A65', 'ISWEB S.p.A.')
    ,@('B65', '01722270665')
    ,@('C65', '3')
    ,@('A66', 'ICCREA Banca SpA')
    ,@('B66', '04774801007')
    ,@('C66', '2')
    ,@('A67', 'KOINE'' SRL')
    ,@('B67', '01934790971')
    ,@('C67', '2')
    ,@('A68', 'Società Almaviva S.p.A.')
    ,@('B68', '08450891000')
    ,@('A69', 'Banco BPM Società per Azioni')
    ,@('B69', '09722490969')
    ,@('A70', 'Ministero dello Sviluppo Economico')
    ,@('B70', '80230390587')
    ,@('A71', 'I.C.A. - Imposte Comunali Affini – s.r.l.')
    ,@('B71', '02478610583')
    ,@('A74', 'Banca Nazionale del Lavoro S.p.A.')
    ,@('B74', '09339391006')
    ,@('A75', 'Noviservice srl')
    ,@('B75', '02789990922')
    ,@('A76', 'Agenzia Italiana del Farmaco - AIFA')
    ,@('B76', '97345810580')
    ,@('A77', 'Engineering Ingegneria Informatica SpA')
    ,@('B77', '00967720285')
    ,@('A78', 'San Marco SPA')
    ,@('B78', '04142440728')
    ,@('A79', 'BANCA MONTE DEI PASCHI DI SIENA')
    ,@('B79', '00884060526')
    ,@('C79', '1')
)

foreach ($change in $changes) {
    $ws.Range($change[0]).Value = $change[1]
}
